# Insert a new data row at row 100 (pushing the existing rows 100-148 down
# to 101-149) and populate it with the new Pina/Tercera record for
# Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 100 and below down by one row.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100.
$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44466
$ws.Range("D100").NumberFormat = $ws.Range("D101").NumberFormat
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100108
$ws.Range("H100").Value = "Tropicales y subtropicales"
$ws.Range("I100").Value = 100108005
$ws.Range("J100").Value = "Piña"
$ws.Range("K100").Value = "Caramelo"
$ws.Range("L100").Value = "Tercera"
$ws.Range("M100").Value = 250
$ws.Range("N100").Value = 20000
$ws.Range("O100").Value = 20000
$ws.Range("P100").Value = 20000
$ws.Range("Q100").Value = "`$/caja 16 unidades"
$ws.Range("R100").Value = "Ecuador"
$ws.Range("S100").Value = 1250
$ws.Range("T100").Value = 16
